$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961298447628"
$ws1.Range("B2").Value = "go_stims-1650996129796759.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961298287687.csv"
$ws1.Range("B4").Value = "go_stims-16509961298287687.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961298447628.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1650996131501465"
$ws2.Range("B2").Value = "ZB-match_1-16509961302047265.csv"
$ws2.Range("B3").Value = "ZB-match_0-16509961299167287.csv"
$ws2.Range("B4").Value = "TB-16509961308767247.csv"
$ws2.Range("B5").Value = "OB-16509961304207244.csv"
$ws2.Range("B6").Value = "ZB-match_3-16509961300047584.csv"
$ws2.Range("B7").Value = "TB-1650996131173467.csv"
$ws2.Range("B8").Value = "TB-16509961314694648.csv"
$ws2.Range("B9").Value = "OB-1650996130380726.csv"
$ws2.Range("B10").Value = "OB-16509961306127632.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650996131501465"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961315494668"
$ws4.Range("B2").Value = "MM_stims-16509961315174663.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996131501465.csv"
$ws4.Range("B4").Value = "MM_stims-16509961315334663.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961315174663.csv"
$ws4.Range("B6").Value = "MM_stims-16509961315494668.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961315334663.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961316296701"
$ws5.Range("B2").Value = "SAT_stims-16509961315815005.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961315494668.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961315976703.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961316137033.csv"
